# Edit script for BAJA_MIPF881205HOCNLR06.docx
# "Cambio en la funcion de reincoorporar empleado"

$d = $word.ActiveDocument

# 1) Mark the "POR OCUPAR OTRO CARGO" checkbox cell in the "MOTIVO DE LA BAJA"
#    table (table 1, row 3, column 2) with an "X ".
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(3)
$cell = $row.Cells.Item(2)
$cell.Range.Text = "X "

# 2) Simple text substitutions throughout the body.
$d.Content.Find.Execute("2A0706A", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2A0508A", 2)

$d.Content.Find.Execute("OFICIAL ADMINISTRATIVO 7A", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "OFICIAL ADMINISTRATIVO 5A", 2)

$d.Content.Find.Execute("SUR 6 MZA.35 LOTE 10 #S/N", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "SUR 6 MZA.35 LOTE 10 FRACC. LOMAS DE NAZARENO", 2)

$d.Content.Find.Execute(" FRACCIONAMIENTO LOMAS DE NAZARENO, SANTA CRUZ XOXOCOTLÁN, OAXACA.", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         " STA CRUZ XOXOCOTLAN", 2)

$d.Content.Find.Execute("undefined", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "71230", 2)

$d.Content.Find.Execute("AUTORIZACIONES Y OTROS PROCEDIMIENTOS", $true, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "CONTROL DE REC. HUMANOS Y SUELDOS APLICADOS", 2)

$d.Content.Find.Execute("1140120000000000310", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "1140020000000000220", 2)

$d.Content.Find.Execute("17 DE OCTUBRE DE 2025", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "3 DE NOVIEMBRE DE 2025", 2)
